# Apply updated computer_score (B) and human_score (C) values for the
# affected rows in the all_agents_score sheet, per the commit:
# "Reward of capture the stag is 4 and human model is better and normalized"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @(2,41,46),
    @(3,24,41),
    @(4,9,14),
    @(5,8,9),
    @(9,10,28),
    @(11,38,47),
    @(17,27,15),
    @(18,12,21),
    @(20,26,23),
    @(22,30,23),
    @(26,5,24),
    @(27,18,24),
    @(28,33,35),
    @(29,28,35),
    @(34,13,26),
    @(36,32,28),
    @(37,14,30),
    @(40,7,26),
    @(41,28,34),
    @(43,29,40),
    @(45,30,17),
    @(46,8,29),
    @(47,13,35),
    @(48,32,41),
    @(49,40,48),
    @(56,9,23),
    @(57,37,51),
    @(58,42,48),
    @(61,19,19),
    @(64,16,30),
    @(74,6,26),
    @(75,14,27),
    @(77,45,30),
    @(79,6,14),
    @(80,39,43),
    @(81,12,19),
    @(82,9,16),
    @(86,34,30),
    @(87,32,48),
    @(88,39,40),
    @(89,30,35),
    @(90,31,30),
    @(91,21,34),
    @(92,39,36),
    @(97,23,38),
    @(98,8,28),
    @(99,27,11),
    @(102,33,29),
    @(103,34,43),
    @(106,30,19),
    @(107,11,22),
    @(110,28,34),
    @(111,39,42),
    @(112,44,48),
    @(113,13,23),
    @(114,10,18),
    @(115,11,25),
    @(118,10,23),
    @(119,13,22),
    @(122,42,46),
    @(123,11,32),
    @(124,14,26),
    @(125,37,42),
    @(127,36,38),
    @(129,12,32),
    @(131,8,26),
    @(133,27,20),
    @(134,28,32),
    @(135,26,20)
)

foreach ($u in $updates) {
    $row = $u[0]
    $newB = $u[1]
    $newC = $u[2]
    $ws.Cells.Item($row, 2).Value = $newB
    $ws.Cells.Item($row, 3).Value = $newC
}
